$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 174, which shifts rows 174:225 down to 175:226
# (the existing row 225 becomes the new row 226, preserving its data/format).
$ws.Rows("174:174").Insert()

# Populate the newly inserted row 174 with the new record's data.
$ws.Range("A174").Value = 11
$ws.Range("B174").Value = "Vega Monumental Concepción"
$ws.Range("C174").Value = "Bíobío"
$ws.Range("D174").Value = 45120
$ws.Range("E174").Value = 8
$ws.Range("F174").Value = 100112043
$ws.Range("G174").Value = "Pepino ensalada"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 100
$ws.Range("K174").Value = 14000
$ws.Range("L174").Value = 15000
$ws.Range("M174").Value = 14500
$ws.Range("N174").Value = "$/caja 60 unidades"
$ws.Range("O174").Value = "Región de Arica y Parinacota"
$ws.Range("P174").Value = 242
$ws.Range("Q174").Value = 60
$ws.Range("R174").Value = "Hortaliza"
